$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the two new data rows (42 and 43) with the new proficiency-level
# value "I" in column B.
$ws.Range("A42").Value = 39
$ws.Range("B42").Value = "I"
$ws.Range("C42").Value = "F"

$ws.Range("A43").Value = 37
$ws.Range("B43").Value = "I"
$ws.Range("C43").Value = "C"

# Apply an AutoFilter over the data range (header row included) filtering
# column A ("Program") down to the discrete value 37 - this both writes the
# <autoFilter> definition and hides the rows that don't match.
$rng = $ws.Range("A1:C42")
$rng.AutoFilter(1, "37", 7)

# Excel records the filter range as a hidden workbook-level defined name
# scoped to the sheet.
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$C`$42")
$fdb.Visible = $false

# Restore the selection to where the author last left it.
$ws.Range("C49").Select()
